$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row: the first block of 10 columns (A:J) that used to
# carry the "_old" suffix now carries "_FV2404" (the older format version),
# the "diff" column (K) is left untouched, and the second block of 10
# columns (L:U) that used to carry the "_new" suffix now carries "_FV2410"
# (the newer format version).
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the data range into a proper Excel Table ("Table1") with an
# autofilter, matching the exported AHB-diff formatting.
$dataRange = $ws.Range("A1:U72")
$table = $ws.ListObjects.Add(1, $dataRange, 0, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# Freeze the header row (top row) so it stays visible while scrolling.
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
